# Apply updated crypto price/volume figures to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value reads as numeric text ("0.999", "40.77", ...)
# must be forced to a text number format first, otherwise Excel would
# silently store them as real numbers instead of strings like the
# original sheet does.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "48.313.07"
$ws.Range("E2").Value = "  +2.11%  "
$ws.Range("D3").Value = "2.528.79"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "323.68"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").Value = "109.46"
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("D7").Value = "0.528"
$ws.Range("E7").Value = "  +0.76%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E9").Value = "  +4.28%  "
$ws.Range("D10").Value = "40.77"
$ws.Range("E10").Value = "  +4.35%  "
$ws.Range("D11").Value = "20.35"
$ws.Range("E11").Value = "  +10.84%  "
$ws.Range("D12").Value = "0.0824"
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("E13").Value = "  +1.14%  "
$ws.Range("D14").Value = "7.29"
$ws.Range("E14").Value = "  +1.56%  "
$ws.Range("D15").Value = "2.925.53"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("D16").Value = "2.527.06"
$ws.Range("E16").Value = "  +0.99%  "
$ws.Range("D17").Value = "0.862"
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("D18").Value = "48.195.43"
$ws.Range("E18").Value = "  +2.03%  "
$ws.Range("D19").Value = "13.34"
$ws.Range("E19").Value = "  +4.12%  "
$ws.Range("D20").Value = "6.65"
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").Value = "0.0₃0948"
$ws.Range("E21").Value = "  +0.66%  "
$ws.Range("E22").Value = "  +0.83%  "
$ws.Range("D23").Value = "72.55"
$ws.Range("E23").Value = "  +2.90%  "
$ws.Range("D24").Value = "269.95"
$ws.Range("E24").Value = "  +8.97%  "
$ws.Range("D25").Value = "2.59"
$ws.Range("E25").Value = "  -0.78%  "
$ws.Range("D26").Value = "26.28"
$ws.Range("E26").Value = "  +0.89%  "
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").Value = "10.21"
$ws.Range("E28").Value = "  +1.36%  "
$ws.Range("E29").Value = "  +6.30%  "
$ws.Range("D30").Value = "35.81"
$ws.Range("E30").Value = "  +1.46%  "
$ws.Range("E31").Value = "  -8.55%  "
$ws.Range("D32").Value = "49.70"
$ws.Range("E32").Value = "  -0.25%  "
$ws.Range("D33").Value = "20.05"
$ws.Range("E33").Value = "  +0.25%  "
$ws.Range("D34").Value = "5.42"
$ws.Range("E34").Value = "  -0.23%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("E37").Value = "  +0.85%  "
$ws.Range("D38").Value = "4.76"
$ws.Range("E38").Value = "  +1.31%  "
$ws.Range("D39").Value = "3.00"
$ws.Range("E39").Value = "  +0.43%  "
$ws.Range("E40").Value = "  +0.20%  "
$ws.Range("D41").Value = "22.45"
$ws.Range("E41").Value = "  +5.20%  "
$ws.Range("E42").Value = "  -2.08%  "
$ws.Range("D43").Value = "118.47"
$ws.Range("E43").Value = "  -2.10%  "
$ws.Range("D44").Value = "0.0300"
$ws.Range("E44").Value = "  +0.47%  "
$ws.Range("D45").Value = "2.014.27"
$ws.Range("E45").Value = "  +1.14%  "
$ws.Range("D46").Value = "3.17"
$ws.Range("E46").Value = "  +3.27%  "
$ws.Range("D47").Value = "1.90"
$ws.Range("E47").Value = "  +6.53%  "
$ws.Range("E48").Value = "  -1.53%  "
$ws.Range("D49").Value = "9.14"
$ws.Range("E49").Value = "  +0.49%  "
$ws.Range("D50").Value = "5.26"
$ws.Range("E50").Value = "  +0.92%  "
$ws.Range("D51").Value = "80.17"
$ws.Range("E51").Value = "  +2.81%  "
